# Auto-generated from the canonical-OOXML diff for Pandaemonium_Profits.xlsx.
# The workbook splits that single logical sheet across 8 worksheets (one per
# crafting job: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Each block below refreshes the
# market-price-derived columns (H:N) of one Leve row with the newly scraped values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 328.77777
$ws.Cells.Item(18, 9).Value = 336
$ws.Cells.Item(18, 11).Value = 336
$ws.Cells.Item(18, 13).Value = -52

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 970.93475
$ws.Cells.Item(129, 9).Value = 397
$ws.Cells.Item(129, 10).Value = 1073.9487
$ws.Cells.Item(129, 11).Value = 1191
$ws.Cells.Item(129, 12).Value = 3221.8461
$ws.Cells.Item(129, 13).Value = 3809
$ws.Cells.Item(129, 14).Value = -13221.8461

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2755.2666
$ws.Cells.Item(2, 9).Value = 1010.1111
$ws.Cells.Item(2, 10).Value = 5373
$ws.Cells.Item(2, 11).Value = 1010.1111
$ws.Cells.Item(2, 12).Value = 5373
$ws.Cells.Item(2, 13).Value = -897.1111
$ws.Cells.Item(2, 14).Value = -5599

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 10592.083
$ws.Cells.Item(61, 9).Value = 11121.223
$ws.Cells.Item(61, 10).Value = 9004.666999999999
$ws.Cells.Item(61, 11).Value = 11121.223
$ws.Cells.Item(61, 12).Value = 9004.666999999999
$ws.Cells.Item(61, 13).Value = -10909.223
$ws.Cells.Item(61, 14).Value = -9428.666999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 5576.5312
$ws.Cells.Item(74, 9).Value = 2337.1177
$ws.Cells.Item(74, 10).Value = 9247.866
$ws.Cells.Item(74, 11).Value = 2337.1177
$ws.Cells.Item(74, 12).Value = 9247.866
$ws.Cells.Item(74, 13).Value = -1463.1177
$ws.Cells.Item(74, 14).Value = -10995.866

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 5576.5312
$ws.Cells.Item(77, 9).Value = 2337.1177
$ws.Cells.Item(77, 10).Value = 9247.866
$ws.Cells.Item(77, 11).Value = 11685.5885
$ws.Cells.Item(77, 12).Value = 46239.33
$ws.Cells.Item(77, 13).Value = -7317.588499999998
$ws.Cells.Item(77, 14).Value = -54975.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3819.889
$ws.Cells.Item(102, 9).Value = 3663.1667
$ws.Cells.Item(102, 10).Value = 4133.3335
$ws.Cells.Item(102, 11).Value = 3663.1667
$ws.Cells.Item(102, 12).Value = 4133.3335
$ws.Cells.Item(102, 13).Value = -2041.1667
$ws.Cells.Item(102, 14).Value = -7377.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2755.2666
$ws.Cells.Item(116, 9).Value = 1010.1111
$ws.Cells.Item(116, 10).Value = 5373
$ws.Cells.Item(116, 11).Value = 1010.1111
$ws.Cells.Item(116, 12).Value = 5373
$ws.Cells.Item(116, 13).Value = 1283.8889
$ws.Cells.Item(116, 14).Value = -9961

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4022.75
$ws.Cells.Item(132, 9).Value = 3656.5
$ws.Cells.Item(132, 11).Value = 10969.5
$ws.Cells.Item(132, 13).Value = -8439.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 10592.083
$ws.Cells.Item(136, 9).Value = 11121.223
$ws.Cells.Item(136, 10).Value = 9004.666999999999
$ws.Cells.Item(136, 11).Value = 33363.669
$ws.Cells.Item(136, 12).Value = 27014.001
$ws.Cells.Item(136, 13).Value = -30813.669
$ws.Cells.Item(136, 14).Value = -32114.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2755.2666
$ws.Cells.Item(3, 9).Value = 1010.1111
$ws.Cells.Item(3, 10).Value = 5373
$ws.Cells.Item(3, 11).Value = 1010.1111
$ws.Cells.Item(3, 12).Value = 5373
$ws.Cells.Item(3, 13).Value = -896.1111
$ws.Cells.Item(3, 14).Value = -5601

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2124.1404
$ws.Cells.Item(86, 9).Value = 1840.8695
$ws.Cells.Item(86, 10).Value = 3308.7273
$ws.Cells.Item(86, 11).Value = 1840.8695
$ws.Cells.Item(86, 12).Value = 3308.7273
$ws.Cells.Item(86, 13).Value = -717.8695
$ws.Cells.Item(86, 14).Value = -5554.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2124.1404
$ws.Cells.Item(89, 9).Value = 1840.8695
$ws.Cells.Item(89, 10).Value = 3308.7273
$ws.Cells.Item(89, 11).Value = 9204.3475
$ws.Cells.Item(89, 12).Value = 16543.6365
$ws.Cells.Item(89, 13).Value = -3588.3475
$ws.Cells.Item(89, 14).Value = -27775.6365

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1441.2
$ws.Cells.Item(94, 9).Value = 1324.7307
$ws.Cells.Item(94, 11).Value = 1324.7307
$ws.Cells.Item(94, 13).Value = -873.7307000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 17933
$ws.Cells.Item(50, 10).Value = 17933
$ws.Cells.Item(50, 12).Value = 17933
$ws.Cells.Item(50, 14).Value = -19183

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 15966
$ws.Cells.Item(51, 10).Value = 15966
$ws.Cells.Item(51, 12).Value = 15966
$ws.Cells.Item(51, 14).Value = -17438

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 23107.643
$ws.Cells.Item(59, 10).Value = 23107.643
$ws.Cells.Item(59, 12).Value = 23107.643
$ws.Cells.Item(59, 14).Value = -25397.643

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 27742.6
$ws.Cells.Item(60, 10).Value = 11345
$ws.Cells.Item(60, 12).Value = 11345
$ws.Cells.Item(60, 14).Value = -12367

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 15966
$ws.Cells.Item(61, 10).Value = 15966
$ws.Cells.Item(61, 12).Value = 15966
$ws.Cells.Item(61, 14).Value = -16662

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3193.6216
$ws.Cells.Item(132, 9).Value = 2639.7334
$ws.Cells.Item(132, 10).Value = 5567.4287
$ws.Cells.Item(132, 11).Value = 7919.2002
$ws.Cells.Item(132, 12).Value = 16702.2861
$ws.Cells.Item(132, 13).Value = -5389.2002
$ws.Cells.Item(132, 14).Value = -21762.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2141.9768
$ws.Cells.Item(134, 9).Value = 1825.7742
$ws.Cells.Item(134, 10).Value = 2958.8333
$ws.Cells.Item(134, 11).Value = 5477.3226
$ws.Cells.Item(134, 12).Value = 8876.499899999999
$ws.Cells.Item(134, 13).Value = -2942.3226
$ws.Cells.Item(134, 14).Value = -13946.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 6948866
$ws.Cells.Item(5, 9).Value = 255.77777
$ws.Cells.Item(5, 10).Value = 27794696
$ws.Cells.Item(5, 11).Value = 767.33331
$ws.Cells.Item(5, 12).Value = 83384088
$ws.Cells.Item(5, 13).Value = -655.33331
$ws.Cells.Item(5, 14).Value = -83384312

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2937.72
$ws.Cells.Item(68, 10).Value = 4438.8965
$ws.Cells.Item(68, 12).Value = 13316.6895
$ws.Cells.Item(68, 14).Value = -14938.6895

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 2937.72
$ws.Cells.Item(71, 10).Value = 4438.8965
$ws.Cells.Item(71, 12).Value = 39950.0685
$ws.Cells.Item(71, 14).Value = -48062.0685

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 6948866
$ws.Cells.Item(135, 9).Value = 255.77777
$ws.Cells.Item(135, 10).Value = 27794696
$ws.Cells.Item(135, 11).Value = 2301.99993
$ws.Cells.Item(135, 12).Value = 250152264
$ws.Cells.Item(135, 13).Value = 233.0000700000001
$ws.Cells.Item(135, 14).Value = -250157334

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1555
$ws.Cells.Item(2, 9).Value = 1555
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1555
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1442
$ws.Cells.Item(2, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 16400
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 16400
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 16400
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).Value = -16624

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1293.3334
$ws.Cells.Item(97, 9).Value = 1340
$ws.Cells.Item(97, 10).Value = 1200
$ws.Cells.Item(97, 11).Value = 1340
$ws.Cells.Item(97, 12).Value = 1200
$ws.Cells.Item(97, 13).Value = -844
$ws.Cells.Item(97, 14).Value = -2192

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3141.1516
$ws.Cells.Item(102, 9).Value = 2823.7368
$ws.Cells.Item(102, 10).Value = 3571.9285
$ws.Cells.Item(102, 11).Value = 2823.7368
$ws.Cells.Item(102, 12).Value = 3571.9285
$ws.Cells.Item(102, 13).Value = -1201.7368
$ws.Cells.Item(102, 14).Value = -6815.9285

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4168.5
$ws.Cells.Item(132, 9).Value = 3116.6667
$ws.Cells.Item(132, 10).Value = 4799.6
$ws.Cells.Item(132, 11).Value = 9350.000100000001
$ws.Cells.Item(132, 12).Value = 14398.8
$ws.Cells.Item(132, 13).Value = -6820.000100000001
$ws.Cells.Item(132, 14).Value = -19458.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 5000
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 5000
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 5000
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(2, 14).Value = -5224

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 2022.375
$ws.Cells.Item(9, 9).Value = 1440
$ws.Cells.Item(9, 10).Value = 2993
$ws.Cells.Item(9, 11).Value = 1440
$ws.Cells.Item(9, 12).Value = 2993
$ws.Cells.Item(9, 13).Value = -1216
$ws.Cells.Item(9, 14).Value = -3441

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1245.6786
$ws.Cells.Item(22, 9).Value = 707.0909
$ws.Cells.Item(22, 10).Value = 1594.1765
$ws.Cells.Item(22, 11).Value = 707.0909
$ws.Cells.Item(22, 12).Value = 1594.1765
$ws.Cells.Item(22, 13).Value = -412.0909
$ws.Cells.Item(22, 14).Value = -2184.1765

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1245.6786
$ws.Cells.Item(27, 9).Value = 707.0909
$ws.Cells.Item(27, 10).Value = 1594.1765
$ws.Cells.Item(27, 11).Value = 707.0909
$ws.Cells.Item(27, 12).Value = 1594.1765
$ws.Cells.Item(27, 13).Value = -600.0909
$ws.Cells.Item(27, 14).Value = -1808.1765

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3307.087
$ws.Cells.Item(40, 9).Value = 3091.9412
$ws.Cells.Item(40, 10).Value = 3916.6667
$ws.Cells.Item(40, 11).Value = 3091.9412
$ws.Cells.Item(40, 12).Value = 3916.6667
$ws.Cells.Item(40, 13).Value = -2955.9412
$ws.Cells.Item(40, 14).Value = -4188.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 400289.6
$ws.Cells.Item(55, 9).Value = 1333438
$ws.Cells.Item(55, 10).Value = 368.85715
$ws.Cells.Item(55, 11).Value = 1333438
$ws.Cells.Item(55, 12).Value = 368.85715
$ws.Cells.Item(55, 13).Value = -1333265
$ws.Cells.Item(55, 14).Value = -714.85715

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 32601.2
$ws.Cells.Item(61, 9).Value = 32601.2
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 32601.2
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -32399.2
$ws.Cells.Item(61, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3621.25
$ws.Cells.Item(68, 10).Value = 4250
$ws.Cells.Item(68, 12).Value = 4250
$ws.Cells.Item(68, 14).Value = -5748

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 3621.25
$ws.Cells.Item(71, 10).Value = 4250
$ws.Cells.Item(71, 12).Value = 21250
$ws.Cells.Item(71, 14).Value = -28738

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2620
$ws.Cells.Item(82, 9).Value = 2000
$ws.Cells.Item(82, 10).Value = 2775
$ws.Cells.Item(82, 11).Value = 2000
$ws.Cells.Item(82, 12).Value = 2775
$ws.Cells.Item(82, 13).Value = -1639
$ws.Cells.Item(82, 14).Value = -3497

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2620
$ws.Cells.Item(85, 9).Value = 2000
$ws.Cells.Item(85, 10).Value = 2775
$ws.Cells.Item(85, 11).Value = 2000
$ws.Cells.Item(85, 12).Value = 2775
$ws.Cells.Item(85, 13).Value = -752
$ws.Cells.Item(85, 14).Value = -5271

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 32601.2
$ws.Cells.Item(113, 9).Value = 32601.2
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 32601.2
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -30431.2
$ws.Cells.Item(113, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2903.6924
$ws.Cells.Item(132, 9).Value = 2036.4166
$ws.Cells.Item(132, 11).Value = 6109.2498
$ws.Cells.Item(132, 13).Value = -3579.2498

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 8235.3125
$ws.Cells.Item(136, 9).Value = 9159.546
$ws.Cells.Item(136, 10).Value = 6202
$ws.Cells.Item(136, 11).Value = 27478.638
$ws.Cells.Item(136, 12).Value = 18606
$ws.Cells.Item(136, 13).Value = -24928.638
$ws.Cells.Item(136, 14).Value = -23706

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 80000000
$ws.Cells.Item(2, 9).Value = 80000000
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 80000000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -79999888
$ws.Cells.Item(2, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1187.7778
$ws.Cells.Item(126, 9).Value = 1201.3334
$ws.Cells.Item(126, 11).Value = 3604.0002
$ws.Cells.Item(126, 13).Value = -1134.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3445
$ws.Cells.Item(132, 9).Value = 3718.6155
$ws.Cells.Item(132, 11).Value = 11155.8465
$ws.Cells.Item(132, 13).Value = -8625.8465
